$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing "Valor"
# column (and its header) from B to C, preserving the B-column style (s="1")
# on row 1 for the new column as well since Excel's InsertBefore copies
# the format of the column to the left.
$ws.Columns("B").Insert()

# New column B: "Variável" header + "Diferença 2021-2012" for every data row
$ws.Range("B1").Value = "Variável"
$ws.Range("B2").Value = "Diferença 2021-2012"
$ws.Range("B3").Value = "Diferença 2021-2012"
$ws.Range("B4").Value = "Diferença 2021-2012"
$ws.Range("B5").Value = "Diferença 2021-2012"
$ws.Range("B6").Value = "Diferença 2021-2012"
$ws.Range("B7").Value = "Diferença 2021-2012"
$ws.Range("B8").Value = "Diferença 2021-2012"
$ws.Range("B9").Value = "Diferença 2021-2012"

# New column D: "Colocação" header + ranking strings for the top 6 rows
$ws.Range("D1").Value = "Colocação"
$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"

$wb.Save()
